$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 13:52"

# Swap province labels: Malaga/Salamanca (rows 19/20) and Cuenca/Avila (rows 40/41)
$ws.Range("A19").Value = "Salamanca"
$ws.Range("A20").Value = "Malaga"
$ws.Range("A40").Value = "Avila"
$ws.Range("A41").Value = "Cuenca"

# Update numeric data for the affected rows (Casos totales, Casos activos, Recuperados, Muertes)
$ws.Range("B18").Value = 2492
$ws.Range("C18").Value = 846
$ws.Range("D18").Value = 1436
$ws.Range("E18").Value = 210

$ws.Range("B19").Value = 2257
$ws.Range("C19").Value = 627
$ws.Range("D19").Value = 1372
$ws.Range("E19").Value = 258

$ws.Range("B20").Value = 2203
$ws.Range("C20").Value = 558
$ws.Range("D20").Value = 1462
$ws.Range("E20").Value = 183

$ws.Range("B25").Value = 1913
$ws.Range("C25").Value = 530
$ws.Range("D25").Value = 1232
$ws.Range("E25").Value = 151

$ws.Range("B29").Value = 1679
$ws.Range("C29").Value = 811
$ws.Range("D29").Value = 613
$ws.Range("E29").Value = 255

$ws.Range("B33").Value = 1232
$ws.Range("C33").Value = 517
$ws.Range("D33").Value = 574
$ws.Range("E33").Value = 141

$ws.Range("B38").Value = 1013
$ws.Range("C38").Value = 242
$ws.Range("D38").Value = 686
$ws.Range("E38").Value = 85

$ws.Range("B40").Value = 952
$ws.Range("C40").Value = 354
$ws.Range("D40").Value = 502
$ws.Range("E40").Value = 96

$ws.Range("B41").Value = 920
$ws.Range("C41").Value = 2532
$ws.Range("D41").Value = 9896
$ws.Range("E41").Value = 135

$ws.Range("B45").Value = 636
$ws.Range("C45").Value = 177
$ws.Range("D45").Value = 409
$ws.Range("E45").Value = 50

$ws.Range("B50").Value = 454
$ws.Range("C50").Value = 162
$ws.Range("D50").Value = 239
$ws.Range("E50").Value = 53
